$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.553524255752563
$ws.Range("B1").Value = 4.379415035247803
$ws.Range("C1").Value = 3.225966930389404
$ws.Range("D1").Value = 1.335103869438171
$ws.Range("E1").Value = 0.9384265542030334
